$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking
# strings (e.g. "1.00", "427.97") are stored as text, matching the
# original inlineStr cell type, not auto-converted to numbers.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "66.292.06"
$ws.Range("E2").Value = "  +4.68%  "
$ws.Range("D3").Value = "3.832.17"
$ws.Range("E3").Value = "  +9.65%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "427.97"
$ws.Range("E5").Value = "  +9.95%  "
$ws.Range("D6").Value = "131.27"
$ws.Range("E6").Value = "  +9.11%  "
$ws.Range("D7").Value = "3.831.63"
$ws.Range("E7").Value = "  +9.71%  "
$ws.Range("D8").Value = "0.612"
$ws.Range("E8").Value = "  +4.67%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "0.731"
$ws.Range("E10").Value = "  +8.64%  "
$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  +4.44%  "
$ws.Range("D12").Value = "0.0000333"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "41.69"
$ws.Range("E13").Value = "  +8.07%  "
$ws.Range("E14").Value = "  +14.57%  "
$ws.Range("D15").Value = "4.447.32"
$ws.Range("E15").Value = "  +10.08%  "
$ws.Range("D16").Value = "15.64"
$ws.Range("E16").Value = "  +24.39%  "
$ws.Range("D17").Value = "3.860.73"
$ws.Range("E17").Value = "  +10.69%  "
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").Value = "20.06"
$ws.Range("E19").Value = "  +7.33%  "
$ws.Range("E20").Value = "  +9.01%  "
$ws.Range("D21").Value = "66.558.79"
$ws.Range("E21").Value = "  +5.04%  "
$ws.Range("D22").Value = "416.51"
$ws.Range("E22").Value = "  +5.99%  "
$ws.Range("D23").Value = "15.04"
$ws.Range("E23").Value = "  +8.43%  "
$ws.Range("D24").Value = "85.12"
$ws.Range("E24").Value = "  +5.29%  "
$ws.Range("D25").Value = "3.11"
$ws.Range("E25").Value = "  +9.07%  "
$ws.Range("D26").Value = "37.35"
$ws.Range("E26").Value = "  +12.36%  "
$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  +16.03%  "
$ws.Range("D28").Value = "3.30"
$ws.Range("E28").Value = "  +11.31%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "5.39"
$ws.Range("E29").Value = "  +4.80%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "9.34"
$ws.Range("E30").Value = "  +37.85%  "
$ws.Range("D31").Value = "13.97"
$ws.Range("E31").Value = "  +18.52%  "
$ws.Range("D32").Value = "721.07"
$ws.Range("E32").Value = "  +8.08%  "
$ws.Range("D33").Value = "0.125"
$ws.Range("E33").Value = "  +14.26%  "
$ws.Range("E34").Value = "  +6.80%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "5.79"
$ws.Range("E36").Value = "  +44.98%  "
$ws.Range("D37").Value = "39.04"
$ws.Range("E37").Value = "  +6.82%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +4.23%  "
$ws.Range("D40").Value = "0.0471"
$ws.Range("E40").Value = "  +8.48%  "
$ws.Range("D41").Value = "0.0₃0728"
$ws.Range("E41").Value = "  +16.54%  "
$ws.Range("E42").Value = "  +7.72%  "
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").Value = "3.28"
$ws.Range("E44").Value = "  +6.99%  "
$ws.Range("E45").Value = "  +5.14%  "
$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").Value = "3.40"
$ws.Range("E46").Value = "  +11.20%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").Value = "0.325"
$ws.Range("E47").Value = "  +18.18%  "
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  +43.80%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.63"
$ws.Range("E49").Value = "  +8.14%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "2.06"
$ws.Range("E50").Value = "  +6.20%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "2.86"
$ws.Range("E51").Value = "  +4.98%  "

# Restore default cell style so no stray number-format style lingers
# on cells that did not have one in the original workbook.
$textRange.Style = "Normal"
